$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "'1.39%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'30.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'11.01%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.185"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.74%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05744"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.19%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.590"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'0.8574"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'5.01%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8734"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.98%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'2.86%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07065"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'2.39%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.02926"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.35%"
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'-0.01%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.001513"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.83%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.04138"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-7.91%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006024"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-93.97%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006014"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-3.26%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D18").Value = "'3.040"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.70%"
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'2.277"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.21%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'2.34%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.03274"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'6.02%"
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = "'1.27%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'3.595"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-3.90%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.42%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001215"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.87%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004510"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'0.48%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'20.31%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'-0.63%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03789"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'4.23%"
$ws.Range("E40").Style = "Normal"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.005714"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-5.84%"
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'2.04%"
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002198"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-16.06%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009655"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'17.59%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005098"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-4.07%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-0.07%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.08891"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-18.40%"
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'-38.73%"
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'-0.07%"
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'-0.07%"
$ws.Range("E50").Style = "Normal"
